# feat(export): ajout d'une colonne pour indiquer les SMS
#
# Inserts a new "Notifications SMS" column (column G) in the "Courriers"
# sheet, shifting the existing "Courriers enregistrés" ... "Passages"
# columns one place to the right, and makes "Courriers" the active sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Courriers")

# Insert a new column at G, shifting existing columns right.
$ws.Columns.Item(7).Insert()

# New header cell, matching the bold/wrap style used by the other header
# cells on row 2 (A2:N2, now B2:O2).
$ws.Cells.Item(2, 7).Value = "Notifications SMS"
$ws.Cells.Item(2, 7).Font.Bold = $true
$ws.Cells.Item(2, 7).WrapText = $true

# Match the column width used for the surrounding header columns.
$ws.Columns.Item(7).ColumnWidth = $ws.Columns.Item(6).ColumnWidth

# Move the selection/active cell as recorded after the edit.
$ws.Range("G3").Select()

# Make "Courriers" the active sheet/tab.
$ws.Activate()
